# The legacy GSC export "Video-Indexing" workbook is refreshed by dropping the
# oldest day (2025-11-14) from the rolling date-indexed table on the "Chart"
# sheet, which shifts every following day's row up by one and shrinks the
# table by a single row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the oldest date (2025-11-14, the first data row under the
# header). Deleting it shifts all subsequent rows (and their
# No video indexed / Video indexed / Impressions values) up by one, so the
# row that used to be "2025-11-15" becomes row 2, etc., and the table ends
# up one row shorter (down to row 88) with no duplicate trailing row.
$ws.Rows.Item(2).Delete()
